$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.295507907867432
$ws.Range("B1").Value = 2.410268783569336
$ws.Range("C1").Value = 2.492666244506836
$ws.Range("D1").Value = 3.244524478912354
$ws.Range("E1").Value = 2.302064895629883
